# Auto-generated PowerShell Excel COM-interop script
# Applies numeric data-refresh updates to the Leve profit-calculation sheets
# (columns H:N = market price / profit figures recomputed by the scheduled runner).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 4200.25
$ws.Range("I4").Value = 3600.3333
$ws.Range("K4").Value = 3600.3333
$ws.Range("M4").Value = -3486.3333
$ws.Range("H40").Value = 2145.3794
$ws.Range("I40").Value = 2039.4
$ws.Range("J40").Value = 2258.9285
$ws.Range("K40").Value = 2039.4
$ws.Range("L40").Value = 2258.9285
$ws.Range("M40").Value = -1864.4
$ws.Range("N40").Value = -2608.9285
$ws.Range("H70").Value = 6967.0586
$ws.Range("I70").Value = 6508.1665
$ws.Range("J70").Value = 7217.364
$ws.Range("K70").Value = 19524.4995
$ws.Range("L70").Value = 21652.092
$ws.Range("M70").Value = -19254.4995
$ws.Range("N70").Value = -22192.092
$ws.Range("H73").Value = 6967.0586
$ws.Range("I73").Value = 6508.1665
$ws.Range("J73").Value = 7217.364
$ws.Range("K73").Value = 19524.4995
$ws.Range("L73").Value = 21652.092
$ws.Range("M73").Value = -18588.4995
$ws.Range("N73").Value = -23524.092
$ws.Range("H86").Value = 5727.7393
$ws.Range("I86").Value = 4998.933
$ws.Range("K86").Value = 4998.933
$ws.Range("M86").Value = -3875.933
$ws.Range("H89").Value = 5727.7393
$ws.Range("I89").Value = 4998.933
$ws.Range("K89").Value = 24994.665
$ws.Range("M89").Value = -19378.665
$ws.Range("H98").Value = 3698.5
$ws.Range("I98").Value = 1118.2
$ws.Range("K98").Value = 1118.2
$ws.Range("M98").Value = 379.8
$ws.Range("H112").Value = 891.2727
$ws.Range("J112").Value = 891.2727
$ws.Range("L112").Value = 2673.8181
$ws.Range("N112").Value = -4889.8181
$ws.Range("H116").Value = 3200.1667
$ws.Range("I116").Value = 3239
$ws.Range("K116").Value = 3239
$ws.Range("M116").Value = 203
$ws.Range("H122").Value = 3698.5
$ws.Range("I122").Value = 1118.2
$ws.Range("K122").Value = 3354.6
$ws.Range("M122").Value = -904.6000000000004
$ws.Range("H138").Value = 3524.8408
$ws.Range("I138").Value = 2940.5833
$ws.Range("J138").Value = 3743.9375
$ws.Range("K138").Value = 8821.749899999999
$ws.Range("L138").Value = 11231.8125
$ws.Range("M138").Value = -3681.749899999999
$ws.Range("N138").Value = -21511.8125
$ws.Range("H141").Value = 4778.75
$ws.Range("I141").Value = 4344.6
$ws.Range("K141").Value = 13033.8
$ws.Range("M141").Value = -7853.800000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4298.0938
$ws.Range("I32").Value = 3784.6667
$ws.Range("K32").Value = 3784.6667
$ws.Range("M32").Value = -3497.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3490.5833
$ws.Range("I99").Value = 2235.4
$ws.Range("K99").Value = 2235.4
$ws.Range("M99").Value = -737.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 14897.588
$ws.Range("I22").Value = 177.27272
$ws.Range("J22").Value = 41884.832
$ws.Range("K22").Value = 177.27272
$ws.Range("L22").Value = 41884.832
$ws.Range("M22").Value = 172.72728
$ws.Range("N22").Value = -42584.832
$ws.Range("H31").Value = 7245.143
$ws.Range("I31").Value = 7531.6665
$ws.Range("K31").Value = 7531.6665
$ws.Range("M31").Value = -7236.6665
$ws.Range("H34").Value = 7245.143
$ws.Range("I34").Value = 7531.6665
$ws.Range("K34").Value = 7531.6665
$ws.Range("M34").Value = -7329.6665
$ws.Range("H62").Value = 205
$ws.Range("I62").Value = 205
$ws.Range("K62").Value = 205
$ws.Range("M62").Value = 419
$ws.Range("H65").Value = 205
$ws.Range("I65").Value = 205
$ws.Range("K65").Value = 1025
$ws.Range("M65").Value = 2095
$ws.Range("H134").Value = 2469.5264
$ws.Range("I134").Value = 2408.4666
$ws.Range("K134").Value = 7225.399800000001
$ws.Range("M134").Value = -4690.399800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 81
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H97").Value = 6945758
$ws.Range("I97").Value = 572.5
$ws.Range("J97").Value = 8930097
$ws.Range("K97").Value = 1717.5
$ws.Range("L97").Value = 26790291
$ws.Range("M97").Value = -1221.5
$ws.Range("N97").Value = -26791283
$ws.Range("H123").Value = 11200.5
$ws.Range("I123").Value = 5314
$ws.Range("J123").Value = 17087
$ws.Range("K123").Value = 15942
$ws.Range("L123").Value = 51261
$ws.Range("M123").Value = -13492
$ws.Range("N123").Value = -56161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2332.9167
$ws.Range("I132").Value = 2249.625
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 6748.875
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -4218.875
$ws.Range("N132").Value = -12558.5
$ws.Range("H136").Value = 26992.875
$ws.Range("J136").Value = 26992.875
$ws.Range("L136").Value = 80978.625
$ws.Range("N136").Value = -86078.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 879.75
$ws.Range("I9").Value = 610
$ws.Range("J9").Value = 1149.5
$ws.Range("K9").Value = 610
$ws.Range("L9").Value = 1149.5
$ws.Range("M9").Value = -386
$ws.Range("N9").Value = -1597.5
$ws.Range("H13").Value = 1050
$ws.Range("I13").Value = 300
$ws.Range("J13").Value = 1800
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 1800
$ws.Range("M13").Value = -160
$ws.Range("N13").Value = -2080
$ws.Range("H132").Value = 3682.75
$ws.Range("J132").Value = 5017
$ws.Range("L132").Value = 15051
$ws.Range("N132").Value = -20111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 29999.5
$ws.Range("I34").Value = 29999.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 29999.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -29796.5
$ws.Range("N34").ClearContents()
$ws.Range("H81").Value = 836830
$ws.Range("I81").Value = 2994.1428
$ws.Range("K81").Value = 5988.2856
$ws.Range("M81").Value = -4927.2856
$ws.Range("H84").Value = 836830
$ws.Range("I84").Value = 2994.1428
$ws.Range("K84").Value = 29941.428
$ws.Range("M84").Value = -24637.428
$ws.Range("H96").Value = 5972.1113
$ws.Range("I96").Value = 6353.8
$ws.Range("K96").Value = 6353.8
$ws.Range("M96").Value = -4980.8
$ws.Range("H107").Value = 223.77777
$ws.Range("I107").Value = 190.57143
$ws.Range("K107").Value = 571.71429
$ws.Range("M107").Value = 1348.28571
